$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.472.69'
$ws.Range('E2').Value = '  -1.56%  '
$ws.Range('D3').Value = '2.458.89'
$ws.Range('E3').Value = '  -1.53%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').Value = '''559.62'
$ws.Range('E5').Value = '  -2.21%  '
$ws.Range('D6').Value = '''163.85'
$ws.Range('E6').Value = '  -1.15%  '
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('E8').Value = '  -1.85%  '
$ws.Range('D9').Value = '2.457.76'
$ws.Range('E9').Value = '  -1.54%  '
$ws.Range('D10').Value = '''0.150'
$ws.Range('E10').Value = '  -5.18%  '
$ws.Range('E11').Value = '  -1.82%  '
$ws.Range('D12').Value = '''0.337'
$ws.Range('E12').Value = '  -5.29%  '
$ws.Range('D13').Value = '''4.81'
$ws.Range('E13').Value = '  -2.20%  '
$ws.Range('D14').Value = '2.907.13'
$ws.Range('E14').Value = '  -1.60%  '
$ws.Range('D15').Value = '68.332.03'
$ws.Range('E15').Value = '  -1.66%  '
$ws.Range('D16').Value = '''0.0000169'
$ws.Range('E16').Value = '  -3.22%  '
$ws.Range('E17').Value = '  -5.36%  '
$ws.Range('D18').Value = '2.447.12'
$ws.Range('E18').Value = '  -2.07%  '
$ws.Range('D19').Value = '''10.93'
$ws.Range('E19').Value = '  -2.10%  '
$ws.Range('B20').Value = 'Uniswap'
$ws.Range('C20').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D20').Value = '''7.18'
$ws.Range('E20').Value = '  -2.94%  '
$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D21').Value = '''341.80'
$ws.Range('E21').Value = '  -1.49%  '
$ws.Range('D22').Value = '''3.77'
$ws.Range('E22').Value = '  -2.88%  '
$ws.Range('E23').Value = '  +0.03%  '
$ws.Range('D24').Value = '''1.87'
$ws.Range('E24').Value = '  -2.88%  '
$ws.Range('D25').Value = '''67.40'
$ws.Range('E25').Value = '  -4.62%  '
$ws.Range('D26').Value = '''1.06'
$ws.Range('E26').Value = '  +5.75%  '
$ws.Range('D27').Value = '''3.70'
$ws.Range('E27').Value = '  -4.99%  '
$ws.Range('E28').Value = '  -1.82%  '
$ws.Range('D29').Value = '''8.11'
$ws.Range('E29').Value = '  -6.05%  '
$ws.Range('D30').Value = '0.0₃0830'
$ws.Range('E30').Value = '  -5.94%  '
$ws.Range('D31').Value = '''7.20'
$ws.Range('E31').Value = '  -8.04%  '
$ws.Range('B32').Value = 'FirstDigitalUSD'
$ws.Range('C32').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D32').Value = '''0.999'
$ws.Range('E32').Value = '  -0.05%  '
$ws.Range('B33').Value = 'Bittensor'
$ws.Range('C33').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D33').Value = '''430.93'
$ws.Range('E33').Value = '  -5.05%  '
$ws.Range('B34').Value = 'Fetch.AI'
$ws.Range('C34').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D34').Value = '''1.16'
$ws.Range('E34').Value = '  -2.52%  '
$ws.Range('B35').Value = 'PancakeSwap'
$ws.Range('C35').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D35').Value = '''1.66'
$ws.Range('E35').Value = '  -3.74%  '
$ws.Range('B36').Value = 'Monero'
$ws.Range('C36').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D36').Value = '''156.70'
$ws.Range('E36').Value = '  +0.50%  '
$ws.Range('B37').Value = 'WhiteBITCoin'
$ws.Range('C37').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D37').Value = '''19.01'
$ws.Range('E37').Value = '  -0.19%  '
$ws.Range('B38').Value = 'USDe'
$ws.Range('C38').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D38').Value = '''1.00'
$ws.Range('E38').Value = '  +0.00%  '
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D39').Value = '''0.110'
$ws.Range('E39').Value = '  -4.52%  '
$ws.Range('B40').Value = 'EthereumClassic'
$ws.Range('C40').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D40').Value = '''17.86'
$ws.Range('E40').Value = '  -2.29%  '
$ws.Range('B41').Value = 'PolygonEcosystemToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D41').Value = '''0.305'
$ws.Range('E41').Value = '  -3.08%  '
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D42').Value = '''4.45'
$ws.Range('E42').Value = '  -4.26%  '
$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D43').Value = '''1.52'
$ws.Range('E43').Value = '  -4.74%  '
$ws.Range('B44').Value = 'ImmutableX'
$ws.Range('C44').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D44').Value = '''1.08'
$ws.Range('E44').Value = '  +1.06%  '
$ws.Range('B45').Value = 'dogwifhat'
$ws.Range('C45').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D45').Value = '''2.07'
$ws.Range('E45').Value = '  -4.89%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').Value = '''134.52'
$ws.Range('E46').Value = '  -3.87%  '
$ws.Range('B47').Value = 'Filecoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D47').Value = '''3.35'
$ws.Range('E47').Value = '  -3.34%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').Value = '''0.0716'
$ws.Range('E48').Value = '  -2.19%  '
$ws.Range('B49').Value = 'ARBITRUM'
$ws.Range('C49').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D49').Value = '''0.481'
$ws.Range('E49').Value = '  -6.50%  '
$ws.Range('B50').Value = 'Mantle'
$ws.Range('C50').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D50').Value = '''0.564'
$ws.Range('E50').Value = '  -1.88%  '
$ws.Range('B51').Value = 'Stellar'
$ws.Range('C51').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D51').Value = '''0.0907'
$ws.Range('E51').Value = '  -2.03%  '
